# Auto-generated edit script: update crypto price/volume table
# (commit: "Updated cryptos list on Mon Feb 27 18:09:26 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '23.379.95'
$ws.Range('E2').Value = '  +0.55%  '
$ws.Range('D3').Value = '1.636.22'
$ws.Range('E3').Value = '  +1.76%  '
$ws.Range('D4').Value = "'1.004"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.47%  '
$ws.Range('D5').Value = "'1.003"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.30%  '
$ws.Range('D6').Value = "'303.64"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.54%  '
$ws.Range('D7').Value = "'0.3786"
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').Value = "'52.33"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.97%  '
$ws.Range('D9').Value = "'0.3627"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.39%  '
$ws.Range('D10').Value = "'1.240"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.93%  '
$ws.Range('D11').Value = "'0.08090"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.72%  '
$ws.Range('D12').Value = "'1.004"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.51%  '
$ws.Range('D13').Value = "'22.79"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.45%  '
$ws.Range('D14').Value = "'6.616"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.06%  '
$ws.Range('D15').Value = "'0.00001248"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.17%  '
$ws.Range('D16').Value = "'7.260"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.39%  '
$ws.Range('D17').Value = '1.634.71'
$ws.Range('E17').Value = '  +1.74%  '
$ws.Range('D18').Value = "'93.94"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.14%  '
$ws.Range('D19').Value = "'0.06938"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.36%  '
$ws.Range('D20').Value = "'18.08"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.63%  '
$ws.Range('D21').Value = "'6.521"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.30%  '
$ws.Range('E22').Value = '  +0.23%  '
$ws.Range('D23').Value = '23.379.97'
$ws.Range('E23').Value = '  +0.53%  '
$ws.Range('D24').Value = "'12.80"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.74%  '
$ws.Range('D25').Value = "'3.246"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +5.03%  '
$ws.Range('D26').Value = "'2.455"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.97%  '
$ws.Range('D27').Value = "'21.14"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.31%  '
$ws.Range('D28').Value = "'149.00"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.13%  '
$ws.Range('D29').Value = "'5.308"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.51%  '
$ws.Range('D30').Value = "'135.28"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.08%  '
$ws.Range('D31').Value = "'2.315"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.84%  '
$ws.Range('D32').Value = '1.817.32'
$ws.Range('E32').Value = '  +1.93%  '
$ws.Range('D33').Value = "'6.845"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.56%  '
$ws.Range('D34').Value = "'10.98"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +6.33%  '
$ws.Range('D35').Value = "'0.9601"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.61%  '
$ws.Range('D36').Value = "'0.02856"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.08%  '
$ws.Range('B37').Value = 'Algorand'
$ws.Range('C37').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D37').Value = "'0.2549"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.43%  '
$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').Value = "'6.224"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.54%  '
$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D39').Value = "'0.08877"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.35%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').Value = "'0.07224"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.29%  '
$ws.Range('D41').Value = "'1.370"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.92%  '
$ws.Range('D42').Value = "'0.7088"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.22%  '
$ws.Range('D43').Value = "'16.39"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.56%  '
$ws.Range('D44').Value = "'12.49"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.07%  '
$ws.Range('D45').Value = "'0.6519"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.25%  '
$ws.Range('D46').Value = "'2.349"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.72%  '
$ws.Range('E47').Value = '  +0.26%  '
$ws.Range('D48').Value = "'3.991"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.50%  '
$ws.Range('D49').Value = "'0.07986"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.25%  '
$ws.Range('D50').Value = "'1.217"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.51%  '
$ws.Range('D51').Value = "'127.22"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.06%  '
